$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: A=1, B=2, C=3, D=4, E=5
# Force the Price column to Text format before writing so values like
# "22.30", "0.9990", "0.00001030" keep their exact digits/trailing zeros
# instead of being auto-coerced to numbers by Excel's smart-entry parsing.

function Set-Price($row, $value) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

function Set-Volume($row, $value) {
    $ws.Cells.Item($row, 5).Value = $value
}

# Rows 2-42: only Price (D) and Volume(1h) (E) values change (row 40's D stays the same)
$updates = @(
    @{ Row=2;  D="29.278.14";   E="  -1.55%  " },
    @{ Row=3;  D="1.900.53";    E="  -2.48%  " },
    @{ Row=4;  D="0.9978";      E="  -0.28%  " },
    @{ Row=5;  D="331.68";      E="  -3.04%  " },
    @{ Row=6;  D="0.9985";      E="  -0.22%  " },
    @{ Row=7;  D="0.4611";      E="  -3.56%  " },
    @{ Row=8;  D="0.4139";      E="  -0.15%  " },
    @{ Row=9;  D="47.69";       E="  -1.25%  " },
    @{ Row=10; D="0.08012";     E="  -2.97%  " },
    @{ Row=11; D="1.009";       E="  -3.04%  " },
    @{ Row=12; D="22.30";       E="  -1.54%  " },
    @{ Row=13; D="1.886.42";    E="  -2.80%  " },
    @{ Row=14; D="5.949";       E="  -4.17%  " },
    @{ Row=15; D="7.110";       E="  -4.21%  " },
    @{ Row=16; D="89.28";       E="  -3.30%  " },
    @{ Row=17; D="0.9990";      E="  -0.26%  " },
    @{ Row=18; D="0.00001030";  E="  -3.17%  " },
    @{ Row=19; D="0.06567";     E="  -1.73%  " },
    @{ Row=20; D="17.64";       E="  -2.25%  " },
    @{ Row=21; D="0.9956";      E="  -0.45%  " },
    @{ Row=22; D="29.212.81";   E="  -1.61%  " },
    @{ Row=23; D="5.478";       E="  -2.27%  " },
    @{ Row=24; D="11.44";       E="  +1.60%  " },
    @{ Row=25; D="2.200";       E="  -3.73%  " },
    @{ Row=26; D="2.148.55";    E="  -1.09%  " },
    @{ Row=27; D="156.52";      E="  -2.86%  " },
    @{ Row=28; D="19.71";       E="  -2.41%  " },
    @{ Row=29; D="2.114";       E="  -3.66%  " },
    @{ Row=30; D="5.618";       E="  -1.24%  " },
    @{ Row=31; D="117.54";      E="  -4.17%  " },
    @{ Row=32; D="1.043";       E="  +1.86%  " },
    @{ Row=33; D="0.09375";     E="  -2.78%  " },
    @{ Row=34; D="1.425";       E="  -3.61%  " },
    @{ Row=35; D="3.524";       E="  -4.31%  " },
    @{ Row=36; D="5.355";       E="  -2.82%  " },
    @{ Row=37; D="0.06077";     E="  -3.85%  " },
    @{ Row=38; D="0.02242";     E="  -3.77%  " },
    @{ Row=39; D="8.438";       E="  -1.82%  " },
    @{ Row=41; D="0.5830";      E="  -4.76%  " },
    @{ Row=42; D="0.9970";      E="  -0.34%  " }
)

foreach ($u in $updates) {
    Set-Price $u.Row $u.D
    Set-Volume $u.Row $u.E
}

# Row 40: only Volume(1h) (E) changes, Price (D) stays "1.176"
Set-Volume 40 "  -1.76%  "

# Row 45: only Price (D) changes; name/link/volume unchanged
Set-Price 45 "1.246"

# Rows 43-51: coin identity/price/volume reshuffle (and a brand-new coin, Quant, at row 51)

# Row 43: Aptos -> Algorand
$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-Price 43 "0.1829"
Set-Volume 43 "  -3.75%  "

# Row 44: Algorand -> Aptos
$ws.Cells.Item(44, 2).Value = "Aptos"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-Price 44 "10.18"
Set-Volume 44 "  -5.31%  "

# Row 46: Cronos -> RenderToken
$ws.Cells.Item(46, 2).Value = "RenderToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-Price 46 "2.322"
Set-Volume 46 "  -2.82%  "

# Row 47: RenderToken -> Cronos
$ws.Cells.Item(47, 2).Value = "Cronos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-Price 47 "0.07536"
Set-Volume 47 "  +1.79%  "

# Row 48: Decentraland -> EnergySwap
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-Price 48 "12.14"
Set-Volume 48 "  -3.72%  "

# Row 49: EnergySwap -> Decentraland
$ws.Cells.Item(49, 2).Value = "Decentraland"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-Price 49 "0.5517"
Set-Volume 49 "  -3.68%  "

# Row 50: PaxosStandard -> NEARProtocol
$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-Price 50 "1.924"
Set-Volume 50 "  -3.92%  "

# Row 51: NEARProtocol -> Quant
$ws.Cells.Item(51, 2).Value = "Quant"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-Price 51 "111.95"
Set-Volume 51 "  -1.94%  "
